$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update RoF values for Rifle (C2) and Pistol (C3)
$ws.Range("C2").Value = 1.5
$ws.Range("C3").Value = 1.5

# Update the active selection to reflect the latest interaction
$ws.Range("L4").Select()
